$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 273.3846
$ws.Range("I9").Value = 278.55554
$ws.Range("J9").Value = 261.75
$ws.Range("K9").Value = 278.55554
$ws.Range("L9").Value = 261.75
$ws.Range("M9").Value = -109.55554
$ws.Range("N9").Value = -599.75
$ws.Range("H33").Value = 7170
$ws.Range("I33").Value = 10267.9
$ws.Range("K33").Value = 10267.9
$ws.Range("M33").Value = -10038.9
$ws.Range("H82").Value = 4409.857
$ws.Range("I82").Value = 1780.4
$ws.Range("J82").Value = 10983.5
$ws.Range("K82").Value = 5341.200000000001
$ws.Range("L82").Value = 32950.5
$ws.Range("M82").Value = -4935.200000000001
$ws.Range("N82").Value = -33762.5
$ws.Range("H85").Value = 4409.857
$ws.Range("I85").Value = 1780.4
$ws.Range("J85").Value = 10983.5
$ws.Range("K85").Value = 5341.200000000001
$ws.Range("L85").Value = 32950.5
$ws.Range("M85").Value = -3937.200000000001
$ws.Range("N85").Value = -35758.5
$ws.Range("H97").Value = 3830.3333
$ws.Range("J97").Value = 3796.4
$ws.Range("L97").Value = 11389.2
$ws.Range("N97").Value = -12381.2
$ws.Range("H101").Value = 399.57144
$ws.Range("J101").Value = 993.5
$ws.Range("L101").Value = 2980.5
$ws.Range("N101").Value = -6224.5
$ws.Range("H132").Value = 6596.407
$ws.Range("I132").Value = 7157.2173
$ws.Range("K132").Value = 21471.6519
$ws.Range("M132").Value = -18941.6519
$ws.Range("H137").Value = 3507.4546
$ws.Range("I137").Value = 1958.4
$ws.Range("J137").Value = 4798.3335
$ws.Range("K137").Value = 5875.200000000001
$ws.Range("L137").Value = 14395.0005
$ws.Range("M137").Value = -3325.200000000001
$ws.Range("N137").Value = -19495.0005
$ws.Range("H138").Value = 2706.9807
$ws.Range("I138").Value = 2276.1614
$ws.Range("K138").Value = 6828.4842
$ws.Range("M138").Value = -1688.4842

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 28302.719
$ws.Range("I32").Value = 28925.166
$ws.Range("K32").Value = 28925.166
$ws.Range("M32").Value = -28638.166
$ws.Range("H132").Value = 80837.695
$ws.Range("I132").Value = 80837.695
$ws.Range("K132").Value = 242513.085
$ws.Range("M132").Value = -239983.085

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 9093918
$ws.Range("I20").Value = 12502674
$ws.Range("J20").Value = 3900
$ws.Range("K20").Value = 12502674
$ws.Range("L20").Value = 3900
$ws.Range("M20").Value = -12502427
$ws.Range("N20").Value = -4394
$ws.Range("H105").Value = 3414.1365
$ws.Range("I105").Value = 3338.6191
$ws.Range("J105").Value = 5000
$ws.Range("K105").Value = 3338.6191
$ws.Range("L105").Value = 5000
$ws.Range("M105").Value = -1591.6191
$ws.Range("N105").Value = -8494
$ws.Range("H134").Value = 3591
$ws.Range("I134").Value = 3306.9167
$ws.Range("J134").Value = 7000
$ws.Range("K134").Value = 9920.750100000001
$ws.Range("L134").Value = 21000
$ws.Range("M134").Value = -7385.750100000001
$ws.Range("N134").Value = -26070

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1919.8
$ws.Range("I16").Value = 1911
$ws.Range("K16").Value = 1911
$ws.Range("M16").Value = -1624
$ws.Range("H23").Value = 3500
$ws.Range("I23").Value = 3500
$ws.Range("K23").Value = 3500
$ws.Range("M23").Value = -3260
$ws.Range("H27").Value = 3500
$ws.Range("I27").Value = 3500
$ws.Range("K27").Value = 3500
$ws.Range("M27").Value = -3308
$ws.Range("H31").Value = 4792.9473
$ws.Range("I31").Value = 2936.5557
$ws.Range("K31").Value = 2936.5557
$ws.Range("M31").Value = -2641.5557
$ws.Range("H34").Value = 4792.9473
$ws.Range("I34").Value = 2936.5557
$ws.Range("K34").Value = 2936.5557
$ws.Range("M34").Value = -2734.5557
$ws.Range("H86").Value = 6798.4443
$ws.Range("I86").Value = 7055.4287
$ws.Range("K86").Value = 7055.4287
$ws.Range("M86").Value = -5932.4287
$ws.Range("H89").Value = 6798.4443
$ws.Range("I89").Value = 7055.4287
$ws.Range("K89").Value = 35277.14350000001
$ws.Range("M89").Value = -29661.14350000001
$ws.Range("H107").Value = 2497.5625
$ws.Range("J107").Value = 3551.2632
$ws.Range("L107").Value = 3551.2632
$ws.Range("N107").Value = -7391.263199999999
$ws.Range("H113").Value = 1919.8
$ws.Range("I113").Value = 1911
$ws.Range("K113").Value = 1911
$ws.Range("M113").Value = 259
$ws.Range("H132").Value = 1772.5
$ws.Range("I132").Value = 1772.5
$ws.Range("K132").Value = 5317.5
$ws.Range("M132").Value = -2787.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H40").Value = 201.38461
$ws.Range("I40").Value = 161.33333
$ws.Range("J40").Value = 235.71428
$ws.Range("K40").Value = 645.33332
$ws.Range("L40").Value = 942.85712
$ws.Range("M40").Value = -576.33332
$ws.Range("N40").Value = -1080.85712
$ws.Range("H41").Value = 5599.8887
$ws.Range("I41").Value = 6607.5386
$ws.Range("J41").Value = 2980
$ws.Range("K41").Value = 19822.6158
$ws.Range("L41").Value = 8940
$ws.Range("M41").Value = -19484.6158
$ws.Range("N41").Value = -9616
$ws.Range("H46").Value = 792.7143
$ws.Range("I46").Value = 820
$ws.Range("J46").Value = 724.5
$ws.Range("K46").Value = 2460
$ws.Range("L46").Value = 2173.5
$ws.Range("M46").Value = -2369
$ws.Range("N46").Value = -2355.5
$ws.Range("H58").Value = 7999
$ws.Range("I58").Value = 4999
$ws.Range("J58").Value = 10999
$ws.Range("K58").Value = 14997
$ws.Range("L58").Value = 32997
$ws.Range("M58").Value = -14869
$ws.Range("N58").Value = -33253
$ws.Range("H97").Value = 432
$ws.Range("I97").Value = 247.25
$ws.Range("J97").Value = 579.8
$ws.Range("K97").Value = 741.75
$ws.Range("L97").Value = 1739.4
$ws.Range("M97").Value = -245.75
$ws.Range("N97").Value = -2731.4
$ws.Range("H107").Value = 937
$ws.Range("J107").Value = 937
$ws.Range("L107").Value = 2811
$ws.Range("N107").Value = -6651
$ws.Range("H127").Value = 8999
$ws.Range("J127").Value = 8999
$ws.Range("L127").Value = 26997
$ws.Range("N127").Value = -36917
$ws.Range("H131").Value = 13228.789
$ws.Range("J131").Value = 17198
$ws.Range("L131").Value = 51594
$ws.Range("N131").Value = -61674
$ws.Range("H132").Value = 0
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 0
$ws.Range("M132").ClearContents()
$ws.Range("N132").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 59352.293
$ws.Range("I107").Value = 83449.164
$ws.Range("J107").Value = 1519.8
$ws.Range("K107").Value = 83449.164
$ws.Range("L107").Value = 1519.8
$ws.Range("M107").Value = -81529.164
$ws.Range("N107").Value = -5359.8
$ws.Range("H122").Value = 3262.6072
$ws.Range("I122").Value = 2321.9412
$ws.Range("K122").Value = 6965.823600000001
$ws.Range("M122").Value = -4515.823600000001
$ws.Range("H132").Value = 113616.11
$ws.Range("I132").Value = 145256.42
$ws.Range("K132").Value = 435769.26
$ws.Range("M132").Value = -433239.26

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 2273.7856
$ws.Range("I16").Value = 2294.923
$ws.Range("J16").Value = 1999
$ws.Range("K16").Value = 2294.923
$ws.Range("L16").Value = 1999
$ws.Range("M16").Value = -2124.923
$ws.Range("N16").Value = -2339
$ws.Range("H46").Value = 10166
$ws.Range("I46").Value = 55750
$ws.Range("J46").Value = 3153.077
$ws.Range("K46").Value = 55750
$ws.Range("L46").Value = 3153.077
$ws.Range("M46").Value = -55562
$ws.Range("N46").Value = -3529.077

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 275
$ws.Range("J2").Value = 300
$ws.Range("L2").Value = 300
$ws.Range("N2").Value = -524
$ws.Range("H107").Value = 1788.8
$ws.Range("I107").Value = 1311.75
$ws.Range("J107").Value = 2334
$ws.Range("K107").Value = 3935.25
$ws.Range("L107").Value = 7002
$ws.Range("M107").Value = -2015.25
$ws.Range("N107").Value = -10842
$ws.Range("H136").Value = 5110.773
$ws.Range("I136").Value = 5588
$ws.Range("J136").Value = 3488.2
$ws.Range("K136").Value = 16764
$ws.Range("L136").Value = 10464.6
$ws.Range("M136").Value = -14214
$ws.Range("N136").Value = -15564.6
